$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-57) need to be sorted ascending by column E (Test Name),
# with blank rows sorted to the bottom. Columns A-L, header in row 1.
$dataRange = $ws.Range("A2:L57")
$keyRange = $ws.Range("E2:E57")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 2
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply()
